$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 23.2600000000002
$ws.Range("H2").Value = [double]"1.43486012875626e-16"
$ws.Range("K2").Value = 55.59143173301403
$ws.Range("L2").Value = "[49.837139098972166, 61.3457243670559]"
$ws.Range("O2").Value = 1.641552918091964
$ws.Range("P2").Value = "[1.540921321580579, 1.7421845146033492]"
$ws.Range("S2").Value = 57.34537162216935
$ws.Range("T2").Value = "[53.57318510774726, 61.11755813659145]"
$ws.Range("W2").Value = 17.18306306306321
$ws.Range("X2").Value = 16.81053053053067
$ws.Range("Y2").Value = 17.55559559559575

# Row 3
$ws.Range("B3").Value = 1
$ws.Range("E3").Value = 22.81000000000013
$ws.Range("H3").Value = [double]"1.43486012875626e-16"
$ws.Range("K3").Value = 50.19912752230111
$ws.Range("L3").Value = "[39.46996771064765, 60.92828733395457]"
$ws.Range("O3").Value = 0.798763297309117
$ws.Range("P3").Value = "[0.5849211547224238, 1.0126054398958102]"
$ws.Range("Q3").Value = [double]"5.224265464676137e-12"
$ws.Range("R3").Value = [double]"5.224265464676137e-12"
$ws.Range("S3").Value = 56.03834397465978
$ws.Range("T3").Value = "[50.62353374463118, 61.45315420468838]"
$ws.Range("W3").Value = 19.91023023023034
$ws.Range("X3").Value = 19.13391391391402
$ws.Range("Y3").Value = 20.68654654654666
